# Update the Novi report worksheet with data through 2021-12-08 (8/12)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 79,4
$arr[0,0] = 44460
$arr[0,1] = 0
$arr[0,2] = 12
$arr[0,3] = 121.6791725816264
$arr[1,0] = 44461
$arr[1,1] = 0
$arr[1,2] = 11
$arr[1,3] = 111.5392415331576
$arr[2,0] = 44462
$arr[2,1] = 1
$arr[2,2] = 10
$arr[2,3] = 101.3993104846887
$arr[3,0] = 44463
$arr[3,1] = 0
$arr[3,2] = 8
$arr[3,3] = 81.11944838775096
$arr[4,0] = 44464
$arr[4,1] = 2
$arr[4,2] = 7
$arr[4,3] = 70.9795173392821
$arr[5,0] = 44465
$arr[5,1] = 0
$arr[5,2] = 6
$arr[5,3] = 60.83958629081322
$arr[6,0] = 44466
$arr[6,1] = 0
$arr[6,2] = 3
$arr[6,3] = 30.41979314540661
$arr[7,0] = 44467
$arr[7,1] = 0
$arr[7,2] = 3
$arr[7,3] = 30.41979314540661
$arr[8,0] = 44468
$arr[8,1] = 0
$arr[8,2] = 3
$arr[8,3] = 30.41979314540661
$arr[9,0] = 44469
$arr[9,1] = 1
$arr[9,2] = 3
$arr[9,3] = 30.41979314540661
$arr[10,0] = 44470
$arr[10,1] = 0
$arr[10,2] = 3
$arr[10,3] = 30.41979314540661
$arr[11,0] = 44471
$arr[11,1] = 0
$arr[11,2] = 1
$arr[11,3] = 10.13993104846887
$arr[12,0] = 44472
$arr[12,1] = 0
$arr[12,2] = 1
$arr[12,3] = 10.13993104846887
$arr[13,0] = 44473
$arr[13,1] = 0
$arr[13,2] = 1
$arr[13,3] = 10.13993104846887
$arr[14,0] = 44474
$arr[14,1] = 0
$arr[14,2] = 1
$arr[14,3] = 10.13993104846887
$arr[15,0] = 44475
$arr[15,1] = 0
$arr[15,2] = 1
$arr[15,3] = 10.13993104846887
$arr[16,0] = 44476
$arr[16,1] = 1
$arr[16,2] = 1
$arr[16,3] = 10.13993104846887
$arr[17,0] = 44477
$arr[17,1] = 0
$arr[17,2] = 1
$arr[17,3] = 10.13993104846887
$arr[18,0] = 44478
$arr[18,1] = 0
$arr[18,2] = 1
$arr[18,3] = 10.13993104846887
$arr[19,0] = 44479
$arr[19,1] = 0
$arr[19,2] = 1
$arr[19,3] = 10.13993104846887
$arr[20,0] = 44480
$arr[20,1] = 0
$arr[20,2] = 1
$arr[20,3] = 10.13993104846887
$arr[21,0] = 44481
$arr[21,1] = 0
$arr[21,2] = 1
$arr[21,3] = 10.13993104846887
$arr[22,0] = 44482
$arr[22,1] = 0
$arr[22,2] = 1
$arr[22,3] = 10.13993104846887
$arr[23,0] = 44483
$arr[23,1] = 0
$arr[23,2] = 0
$arr[23,3] = 0
$arr[24,0] = 44484
$arr[24,1] = 0
$arr[24,2] = 0
$arr[24,3] = 0
$arr[25,0] = 44485
$arr[25,1] = 0
$arr[25,2] = 0
$arr[25,3] = 0
$arr[26,0] = 44486
$arr[26,1] = 0
$arr[26,2] = 0
$arr[26,3] = 0
$arr[27,0] = 44487
$arr[27,1] = 0
$arr[27,2] = 0
$arr[27,3] = 0
$arr[28,0] = 44488
$arr[28,1] = 0
$arr[28,2] = 0
$arr[28,3] = 0
$arr[29,0] = 44489
$arr[29,1] = 0
$arr[29,2] = 0
$arr[29,3] = 0
$arr[30,0] = 44490
$arr[30,1] = 0
$arr[30,2] = 0
$arr[30,3] = 0
$arr[31,0] = 44491
$arr[31,1] = 0
$arr[31,2] = 0
$arr[31,3] = 0
$arr[32,0] = 44492
$arr[32,1] = 0
$arr[32,2] = 0
$arr[32,3] = 0
$arr[33,0] = 44493
$arr[33,1] = 0
$arr[33,2] = 0
$arr[33,3] = 0
$arr[34,0] = 44494
$arr[34,1] = 0
$arr[34,2] = 0
$arr[34,3] = 0
$arr[35,0] = 44495
$arr[35,1] = 0
$arr[35,2] = 0
$arr[35,3] = 0
$arr[36,0] = 44496
$arr[36,1] = 0
$arr[36,2] = 0
$arr[36,3] = 0
$arr[37,0] = 44497
$arr[37,1] = 0
$arr[37,2] = 0
$arr[37,3] = 0
$arr[38,0] = 44498
$arr[38,1] = 0
$arr[38,2] = 0
$arr[38,3] = 0
$arr[39,0] = 44499
$arr[39,1] = 0
$arr[39,2] = 0
$arr[39,3] = 0
$arr[40,0] = 44500
$arr[40,1] = 0
$arr[40,2] = 0
$arr[40,3] = 0
$arr[41,0] = 44501
$arr[41,1] = 0
$arr[41,2] = 0
$arr[41,3] = 0
$arr[42,0] = 44502
$arr[42,1] = 0
$arr[42,2] = 0
$arr[42,3] = 0
$arr[43,0] = 44503
$arr[43,1] = 0
$arr[43,2] = 0
$arr[43,3] = 0
$arr[44,0] = 44504
$arr[44,1] = 0
$arr[44,2] = 0
$arr[44,3] = 0
$arr[45,0] = 44505
$arr[45,1] = 0
$arr[45,2] = 0
$arr[45,3] = 0
$arr[46,0] = 44506
$arr[46,1] = 0
$arr[46,2] = 0
$arr[46,3] = 0
$arr[47,0] = 44507
$arr[47,1] = 2
$arr[47,2] = 2
$arr[47,3] = 20.27986209693774
$arr[48,0] = 44508
$arr[48,1] = 0
$arr[48,2] = 2
$arr[48,3] = 20.27986209693774
$arr[49,0] = 44509
$arr[49,1] = 0
$arr[49,2] = 2
$arr[49,3] = 20.27986209693774
$arr[50,0] = 44510
$arr[50,1] = 0
$arr[50,2] = 2
$arr[50,3] = 20.27986209693774
$arr[51,0] = 44511
$arr[51,1] = 0
$arr[51,2] = 2
$arr[51,3] = 20.27986209693774
$arr[52,0] = 44512
$arr[52,1] = 0
$arr[52,2] = 2
$arr[52,3] = 20.27986209693774
$arr[53,0] = 44513
$arr[53,1] = 0
$arr[53,2] = 2
$arr[53,3] = 20.27986209693774
$arr[54,0] = 44514
$arr[54,1] = 0
$arr[54,2] = 0
$arr[54,3] = 0
$arr[55,0] = 44515
$arr[55,1] = 2
$arr[55,2] = 2
$arr[55,3] = 20.27986209693774
$arr[56,0] = 44516
$arr[56,1] = 10
$arr[56,2] = 12
$arr[56,3] = 121.6791725816264
$arr[57,0] = 44517
$arr[57,1] = 0
$arr[57,2] = 12
$arr[57,3] = 121.6791725816264
$arr[58,0] = 44518
$arr[58,1] = 0
$arr[58,2] = 12
$arr[58,3] = 121.6791725816264
$arr[59,0] = 44519
$arr[59,1] = 0
$arr[59,2] = 12
$arr[59,3] = 121.6791725816264
$arr[60,0] = 44520
$arr[60,1] = 1
$arr[60,2] = 13
$arr[60,3] = 131.8191036300953
$arr[61,0] = 44521
$arr[61,1] = 0
$arr[61,2] = 13
$arr[61,3] = 131.8191036300953
$arr[62,0] = 44522
$arr[62,1] = 0
$arr[62,2] = 11
$arr[62,3] = 111.5392415331576
$arr[63,0] = 44523
$arr[63,1] = 0
$arr[63,2] = 1
$arr[63,3] = 10.13993104846887
$arr[64,0] = 44524
$arr[64,1] = 9
$arr[64,2] = 10
$arr[64,3] = 101.3993104846887
$arr[65,0] = 44525
$arr[65,1] = 0
$arr[65,2] = 10
$arr[65,3] = 101.3993104846887
$arr[66,0] = 44526
$arr[66,1] = 0
$arr[66,2] = 10
$arr[66,3] = 101.3993104846887
$arr[67,0] = 44527
$arr[67,1] = 3
$arr[67,2] = 12
$arr[67,3] = 121.6791725816264
$arr[68,0] = 44528
$arr[68,1] = 0
$arr[68,2] = 12
$arr[68,3] = 121.6791725816264
$arr[69,0] = 44529
$arr[69,1] = 1
$arr[69,2] = 13
$arr[69,3] = 131.8191036300953
$arr[70,0] = 44530
$arr[70,1] = 3
$arr[70,2] = 16
$arr[70,3] = 162.2388967755019
$arr[71,0] = 44531
$arr[71,1] = 0
$arr[71,2] = 7
$arr[71,3] = 70.9795173392821
$arr[72,0] = 44532
$arr[72,1] = 5
$arr[72,2] = 12
$arr[72,3] = 121.6791725816264
$arr[73,0] = 44533
$arr[73,1] = 2
$arr[73,2] = 14
$arr[73,3] = 141.9590346785642
$arr[74,0] = 44534
$arr[74,1] = 3
$arr[74,2] = 14
$arr[74,3] = 141.9590346785642
$arr[75,0] = 44535
$arr[75,1] = 2
$arr[75,2] = 16
$arr[75,3] = 162.2388967755019
$arr[76,0] = 44536
$arr[76,1] = 0
$arr[76,2] = 15
$arr[76,3] = 152.0989657270331
$arr[77,0] = 44537
$arr[77,1] = 2
$arr[77,2] = 14
$arr[77,3] = 141.9590346785642
$arr[78,0] = 44538
$arr[78,1] = 1
$arr[78,2] = 15
$arr[78,3] = 152.0989657270331

$ws.Range("A386:D464").Value = $arr

# Copy the date column formatting (style) from the last pre-existing row
$ws.Range("A385").Copy()
$ws.Range("A386:A464").PasteSpecial(-4122)
